$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Find the last used row in the sheet
$lastRow = $ws.UsedRange.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # Column G = "Recorded By"
    $val = $cell.Value2

    if ($val -eq "dnasr281@gmail.com, System") {
        $cell.Value2 = "System, dnasr281@gmail.com"
    }
    elseif ($val -eq "system, backup@backdoor.com, System") {
        $cell.Value2 = "backup@backdoor.com, System, system"
    }
}
